# ---------------------------------------------------------------------------
# Applies the "added changes to match IMB" edit to web-stranica.docx
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$enDash = [char]0x2013

# ---------------------------------------------------------------------------
# 1. Merge the split run around the stray mid-word "_GoBack" bookmark in the
#    "-uprava ... -sistem za navodnjavanje" line into one contiguous run, and
#    drop the now-redundant bookmark that used to sit in the middle of the
#    word "sistem" (it gets relocated to the very end of the document below).
# ---------------------------------------------------------------------------
$oldLine = "-uprava               -sadnice maline               -sistem za navodnjavanje " + $enDash + "agrocentar vrbanja"
$d.Content.Find.Execute($oldLine, $true, $false, $false, $false, $false, $true, 1, $false, $oldLine, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Seed the numbering part + the "List Paragraph" style definition the same
#    way Word does the first time a numbered list is applied. Assigning the
#    style by name first (before touching ListFormat) avoids the host
#    mis-binding the live style object to "Normal"; only then do we tune the
#    style's properties and wire up the default numbered-list formatting.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$scratch = $d.Paragraphs.Item($lastIndex)
$scratch.Style = "List Paragraph"

$listStyle = $d.Styles("List Paragraph")
$listStyle.Priority = 34
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

$scratch.Range.ListFormat.ApplyNumberDefault()

# ---------------------------------------------------------------------------
# 3. Replace the trailing (scratch) paragraph with: a manual page break
#    paragraph, a "Potrebno:" paragraph, and two numbered ("List
#    Paragraph"/numId 1) paragraphs -- the second one now carrying the
#    relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$newBlockXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Potrebno:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Novi logo, da je rezolucija bolja i da su boje drugacije (u vise boja pozeljno) zuto-crveni ne ide sa ovom temom</w:t></w:r><w:r><w:t xml:space="preserve"> (dimenzije pozeljne: )</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Imate li Facebook stranicu, Instagram</w:t></w:r><w:r><w:t>, Twiter..? Ako nema valjalo bi napraviti i tako se oglasavati javnosti a ostatak materijala drzati na web stranici $enDash to je najbolja praksa danas.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$scratch.Range.InsertXML($newBlockXml)
